# Notification services and file upload service skeleton
# Adds a new "File Upload" service row (row 29) to the ServicesList sheet,
# mirroring the layout/formulas of the existing rows, then moves the
# selection to the newly added formula cell (M29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (borders/center alignment) of the last data row (28)
# down onto the new row (29) before filling in values, so the new row
# matches the existing table style (style index 4 / 5).
$ws.Range("B28:L28").Copy()
$ws.Range("B29:L29").PasteSpecial(-4122)

# --- Row 29 data -----------------------------------------------------
$ws.Range("B29").Value = "File Upload"
$ws.Range("C29").Value = "File Upload"
$ws.Range("D29").Value = "WS-FL-01"
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = $false
$ws.Range("G29").Value = "booking/file"
$ws.Range("H29").Value = "/fileupload"
$ws.Range("I29").Value = "POST"

# --- Row 29 helper formulas (INSERT statement / mapping annotation / ---
# --- ServiceInfo annotation), matching columns M/N/O of earlier rows ---
$ws.Range("M29").Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D29,"'',''CONNON_CONFIG'', 0, ''",C29,"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
$ws.Range("N29").Formula = '=_xlfn.CONCAT(IF(I29="GET","@GetMapping(",IF(I29="POST","@PostMapping(",IF(I29="DELETE","@DeleteMapping(",IF(I29="PUT","@PutMapping(","")))),CHAR(34),H29,CHAR(34),")")'
$ws.Range("O29").Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D29,,CHAR(34),", serviceName = ",CHAR(34),C29,CHAR(34), ", queryId = ",CHAR(34),E29,CHAR(34),", logActivity =",F29,")")'

# Move the selection to the newly-entered formula cell, as in the source
# edit (the user ended up with M29 selected after typing the formula).
[void]$ws.Range("M29").Select()
